$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.127.46'
$ws.Range("E2").Value = '  -0.72%  '
$ws.Range("D3").Value = '3.318.23'
$ws.Range("E3").Value = '  -0.50%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '181.94'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.82%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.652'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.10%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '3.315.15'
$ws.Range("E9").Value = '  -0.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.126'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.80'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.45%  '
$ws.Range("E12").Value = '  -0.54%  '
$ws.Range("D13").Value = '3.890.66'
$ws.Range("E13").Value = '  -0.60%  '
$ws.Range("E14").Value = '  -2.82%  '
$ws.Range("D15").Value = '66.164.11'
$ws.Range("E15").Value = '  -0.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.18'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.04%  '
$ws.Range("E17").Value = '  -1.07%  '
$ws.Range("D18").Value = '3.268.38'
$ws.Range("E18").Value = '  -2.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '425.82'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.54%  '
$ws.Range("E20").Value = '  -2.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.13'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.92%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.38'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.73%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.67'
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("E25").Value = '  +0.40%  '
$ws.Range("D26").Value = '3.460.74'
$ws.Range("E26").Value = '  -0.70%  '
$ws.Range("E27").Value = '  -0.83%  '
$ws.Range("E28").Value = '  +5.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000114'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.90%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.88'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.52%  '
$ws.Range("E31").Value = '  +0.14%  '
$ws.Range("E32").Value = '  -2.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '22.37'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.89%  '
$ws.Range("E34").Value = '  +0.10%  '
$ws.Range("E35").Value = '  -1.57%  '
$ws.Range("E36").Value = '  -2.97%  '
$ws.Range("E37").Value = '  -4.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '160.80'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.69%  '
$ws.Range("E39").Value = '  -2.96%  '
$ws.Range("D40").Value = '2.869.83'
$ws.Range("E40").Value = '  +1.54%  '
$ws.Range("E41").Value = '  +0.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '26.39'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.22%  '
$ws.Range("E44").Value = '  -2.45%  '
$ws.Range("E45").Value = '  -1.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0660'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.88%  '
$ws.Range("E47").Value = '  -4.81%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.29'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '23.11'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '312.88'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.72%  '
$ws.Range("E51").Value = '  -0.85%  '
